$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 24 to make room for the new "LEO" entry,
# shifting existing rows 24-51 down to 25-52.
$ws.Rows("24:24").Insert()

# Re-write all data rows (2-52) for columns A-E to match the refreshed data set.
$ws.Cells.Item(2, 1).Value = 0
$ws.Cells.Item(2, 2).Value = 'Bitcoin'
$ws.Cells.Item(2, 3).Value = 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'
$ws.Cells.Item(2, 4).Value = '56.408.39'
$ws.Cells.Item(2, 5).Value = '  -3.59%  '
$ws.Cells.Item(3, 1).Value = 1
$ws.Cells.Item(3, 2).Value = 'Ethereum'
$ws.Cells.Item(3, 3).Value = 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'
$ws.Cells.Item(3, 4).Value = '2.379.50'
$ws.Cells.Item(3, 5).Value = '  -3.98%  '
$ws.Cells.Item(4, 1).Value = 2
$ws.Cells.Item(4, 2).Value = 'TetherUSD'
$ws.Cells.Item(4, 3).Value = 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'
$ws.Cells.Item(4, 4).Value = '1.00'
$ws.Cells.Item(4, 5).Value = '  -0.03%  '
$ws.Cells.Item(5, 1).Value = 3
$ws.Cells.Item(5, 2).Value = 'BNB'
$ws.Cells.Item(5, 3).Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Cells.Item(5, 4).Value = '499.50'
$ws.Cells.Item(5, 5).Value = '  -5.61%  '
$ws.Cells.Item(6, 1).Value = 4
$ws.Cells.Item(6, 2).Value = 'Solana'
$ws.Cells.Item(6, 3).Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Cells.Item(6, 4).Value = '129.11'
$ws.Cells.Item(6, 5).Value = '  -2.26%  '
$ws.Cells.Item(7, 1).Value = 5
$ws.Cells.Item(7, 2).Value = 'USDC'
$ws.Cells.Item(7, 3).Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Cells.Item(7, 4).Value = '0.996'
$ws.Cells.Item(7, 5).Value = '  -0.74%  '
$ws.Cells.Item(8, 1).Value = 6
$ws.Cells.Item(8, 2).Value = 'XRP'
$ws.Cells.Item(8, 3).Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Cells.Item(8, 4).Value = '0.551'
$ws.Cells.Item(8, 5).Value = '  -2.17%  '
$ws.Cells.Item(9, 1).Value = 7
$ws.Cells.Item(9, 2).Value = 'LidoStakedEther'
$ws.Cells.Item(9, 3).Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Cells.Item(9, 4).Value = '2.405.13'
$ws.Cells.Item(9, 5).Value = '  -2.94%  '
$ws.Cells.Item(10, 1).Value = 8
$ws.Cells.Item(10, 2).Value = 'Dogecoin'
$ws.Cells.Item(10, 3).Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Cells.Item(10, 4).Value = '0.0957'
$ws.Cells.Item(10, 5).Value = '  -3.17%  '
$ws.Cells.Item(11, 1).Value = 9
$ws.Cells.Item(11, 2).Value = 'TRON'
$ws.Cells.Item(11, 3).Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Cells.Item(11, 4).Value = '0.150'
$ws.Cells.Item(11, 5).Value = '  -1.02%  '
$ws.Cells.Item(12, 1).Value = 10
$ws.Cells.Item(12, 2).Value = 'Toncoin'
$ws.Cells.Item(12, 3).Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Cells.Item(12, 4).Value = '4.67'
$ws.Cells.Item(12, 5).Value = '  -7.12%  '
$ws.Cells.Item(13, 1).Value = 11
$ws.Cells.Item(13, 2).Value = 'Cardano'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Cells.Item(13, 4).Value = '0.317'
$ws.Cells.Item(13, 5).Value = '  -2.28%  '
$ws.Cells.Item(14, 1).Value = 12
$ws.Cells.Item(14, 2).Value = 'WrappedliquidstakedEther2.0'
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Cells.Item(14, 4).Value = '2.803.06'
$ws.Cells.Item(14, 5).Value = '  -4.14%  '
$ws.Cells.Item(15, 1).Value = 13
$ws.Cells.Item(15, 2).Value = 'WrappedBTC'
$ws.Cells.Item(15, 3).Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Cells.Item(15, 4).Value = '56.275.09'
$ws.Cells.Item(15, 5).Value = '  -3.61%  '
$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(16, 2).Value = 'Avalanche'
$ws.Cells.Item(16, 3).Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Cells.Item(16, 4).Value = '21.51'
$ws.Cells.Item(16, 5).Value = '  -2.33%  '
$ws.Cells.Item(17, 1).Value = 15
$ws.Cells.Item(17, 2).Value = 'ShibaInu'
$ws.Cells.Item(17, 3).Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Cells.Item(17, 4).Value = '0.0000131'
$ws.Cells.Item(17, 5).Value = '  -2.08%  '
$ws.Cells.Item(18, 1).Value = 16
$ws.Cells.Item(18, 2).Value = 'WrappedEther'
$ws.Cells.Item(18, 3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Cells.Item(18, 4).Value = '2.412.57'
$ws.Cells.Item(18, 5).Value = '  -3.10%  '
$ws.Cells.Item(19, 1).Value = 17
$ws.Cells.Item(19, 2).Value = 'Chainlink'
$ws.Cells.Item(19, 3).Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Cells.Item(19, 4).Value = '10.10'
$ws.Cells.Item(19, 5).Value = '  -3.64%  '
$ws.Cells.Item(20, 1).Value = 18
$ws.Cells.Item(20, 2).Value = 'BitcoinCash'
$ws.Cells.Item(20, 3).Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Cells.Item(20, 4).Value = '309.44'
$ws.Cells.Item(20, 5).Value = '  -2.74%  '
$ws.Cells.Item(21, 1).Value = 19
$ws.Cells.Item(21, 2).Value = 'Polkadot'
$ws.Cells.Item(21, 3).Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Cells.Item(21, 4).Value = '4.02'
$ws.Cells.Item(21, 5).Value = '  -4.14%  '
$ws.Cells.Item(22, 1).Value = 20
$ws.Cells.Item(22, 2).Value = 'Uniswap'
$ws.Cells.Item(22, 3).Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Cells.Item(22, 4).Value = '6.27'
$ws.Cells.Item(22, 5).Value = '  +2.20%  '
$ws.Cells.Item(23, 1).Value = 21
$ws.Cells.Item(23, 2).Value = 'Dai'
$ws.Cells.Item(23, 3).Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Cells.Item(23, 4).Value = '0.998'
$ws.Cells.Item(23, 5).Value = '  -0.24%  '
$ws.Cells.Item(24, 1).Value = 22
$ws.Cells.Item(24, 2).Value = 'LEO'
$ws.Cells.Item(24, 3).Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Cells.Item(24, 4).Value = '5.55'
$ws.Cells.Item(24, 5).Value = '  -5.34%  '
$ws.Cells.Item(25, 1).Value = 23
$ws.Cells.Item(25, 2).Value = 'Litecoin'
$ws.Cells.Item(25, 3).Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Cells.Item(25, 4).Value = '64.94'
$ws.Cells.Item(25, 5).Value = '  -0.44%  '
$ws.Cells.Item(26, 1).Value = 24
$ws.Cells.Item(26, 2).Value = 'Binance-PegBSC-USD'
$ws.Cells.Item(26, 3).Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Cells.Item(26, 4).Value = '0.996'
$ws.Cells.Item(26, 5).Value = '  -0.96%  '
$ws.Cells.Item(27, 1).Value = 25
$ws.Cells.Item(27, 2).Value = 'WrappedeETH'
$ws.Cells.Item(27, 3).Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Cells.Item(27, 4).Value = '2.488.24'
$ws.Cells.Item(27, 5).Value = '  -4.98%  '
$ws.Cells.Item(28, 1).Value = 26
$ws.Cells.Item(28, 2).Value = 'Polygon'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Cells.Item(28, 4).Value = '0.376'
$ws.Cells.Item(28, 5).Value = '  -6.74%  '
$ws.Cells.Item(29, 1).Value = 27
$ws.Cells.Item(29, 2).Value = 'Kaspa'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Cells.Item(29, 4).Value = '0.150'
$ws.Cells.Item(29, 5).Value = '  -5.46%  '
$ws.Cells.Item(30, 1).Value = 28
$ws.Cells.Item(30, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Cells.Item(30, 4).Value = '7.27'
$ws.Cells.Item(30, 5).Value = '  -0.74%  '
$ws.Cells.Item(31, 1).Value = 29
$ws.Cells.Item(31, 2).Value = 'Monero'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(31, 4).Value = '172.37'
$ws.Cells.Item(31, 5).Value = '  -1.56%  '
$ws.Cells.Item(32, 1).Value = 30
$ws.Cells.Item(32, 2).Value = 'PEPE'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Cells.Item(32, 4).Value = '0.0₃0713'
$ws.Cells.Item(32, 5).Value = '  -4.07%  '
$ws.Cells.Item(33, 1).Value = 31
$ws.Cells.Item(33, 2).Value = 'PancakeSwap'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Cells.Item(33, 4).Value = '1.65'
$ws.Cells.Item(33, 5).Value = '  -3.30%  '
$ws.Cells.Item(34, 1).Value = 32
$ws.Cells.Item(34, 2).Value = 'Aptos'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Cells.Item(34, 4).Value = '6.12'
$ws.Cells.Item(34, 5).Value = '  -0.46%  '
$ws.Cells.Item(35, 1).Value = 33
$ws.Cells.Item(35, 2).Value = 'USDe'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Cells.Item(35, 4).Value = '0.998'
$ws.Cells.Item(35, 5).Value = '  -0.13%  '
$ws.Cells.Item(36, 1).Value = 34
$ws.Cells.Item(36, 2).Value = 'Fetch.AI'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Cells.Item(36, 4).Value = '1.10'
$ws.Cells.Item(36, 5).Value = '  -6.30%  '
$ws.Cells.Item(37, 1).Value = 35
$ws.Cells.Item(37, 2).Value = 'FirstDigitalUSD'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Cells.Item(37, 4).Value = '0.995'
$ws.Cells.Item(37, 5).Value = '  -0.26%  '
$ws.Cells.Item(38, 1).Value = 36
$ws.Cells.Item(38, 2).Value = 'EthereumClassic'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Cells.Item(38, 4).Value = '17.84'
$ws.Cells.Item(38, 5).Value = '  -0.54%  '
$ws.Cells.Item(39, 1).Value = 37
$ws.Cells.Item(39, 2).Value = 'ImmutableX'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(39, 4).Value = '1.22'
$ws.Cells.Item(39, 5).Value = '  +1.39%  '
$ws.Cells.Item(40, 1).Value = 38
$ws.Cells.Item(40, 2).Value = 'NEARProtocol'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Cells.Item(40, 4).Value = '3.78'
$ws.Cells.Item(40, 5).Value = '  -1.77%  '
$ws.Cells.Item(41, 1).Value = 39
$ws.Cells.Item(41, 2).Value = 'OKB'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Cells.Item(41, 4).Value = '35.87'
$ws.Cells.Item(41, 5).Value = '  -1.30%  '
$ws.Cells.Item(42, 1).Value = 40
$ws.Cells.Item(42, 2).Value = 'SuiNetwork'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Cells.Item(42, 4).Value = '0.789'
$ws.Cells.Item(42, 5).Value = '  -2.54%  '
$ws.Cells.Item(43, 1).Value = 41
$ws.Cells.Item(43, 2).Value = 'Stacks'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Cells.Item(43, 4).Value = '1.42'
$ws.Cells.Item(43, 5).Value = '  -3.76%  '
$ws.Cells.Item(44, 1).Value = 42
$ws.Cells.Item(44, 2).Value = 'Aave'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Cells.Item(44, 4).Value = '129.58'
$ws.Cells.Item(44, 5).Value = '  +0.57%  '
$ws.Cells.Item(45, 1).Value = 43
$ws.Cells.Item(45, 2).Value = 'Filecoin'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Cells.Item(45, 4).Value = '3.33'
$ws.Cells.Item(45, 5).Value = '  -3.09%  '
$ws.Cells.Item(46, 1).Value = 44
$ws.Cells.Item(46, 2).Value = 'RenderToken'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(46, 4).Value = '4.84'
$ws.Cells.Item(46, 5).Value = '  -1.39%  '
$ws.Cells.Item(47, 1).Value = 45
$ws.Cells.Item(47, 2).Value = 'Bittensor'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Cells.Item(47, 4).Value = '251.93'
$ws.Cells.Item(47, 5).Value = '  -6.82%  '
$ws.Cells.Item(48, 1).Value = 46
$ws.Cells.Item(48, 2).Value = 'Mantle'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Cells.Item(48, 4).Value = '0.565'
$ws.Cells.Item(48, 5).Value = '  -3.54%  '
$ws.Cells.Item(49, 1).Value = 47
$ws.Cells.Item(49, 2).Value = 'Stellar'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Cells.Item(49, 4).Value = '0.0898'
$ws.Cells.Item(49, 5).Value = '  -3.80%  '
$ws.Cells.Item(50, 1).Value = 48
$ws.Cells.Item(50, 2).Value = 'Hedera'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(50, 4).Value = '0.0484'
$ws.Cells.Item(50, 5).Value = '  -2.88%  '
$ws.Cells.Item(51, 1).Value = 49
$ws.Cells.Item(51, 2).Value = 'EnergySwap'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(51, 4).Value = '16.83'
$ws.Cells.Item(51, 5).Value = '  -1.16%  '
$ws.Cells.Item(52, 1).Value = 50
$ws.Cells.Item(52, 2).Value = 'VeChain'
$ws.Cells.Item(52, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(52, 4).Value = '0.0206'
$ws.Cells.Item(52, 5).Value = '  -4.74%  '
